# This presentation ships two theme parts:
#   ppt/theme/theme1.xml  -> linked from the slide master (currently "Integral" / Red Violet)
#   ppt/theme/theme2.xml  -> linked from the notes master  (currently "Office Theme")
#
# The target edit swaps the two themes' contents: the slide master's theme
# becomes the stock "Office Theme" palette and the notes master's theme
# becomes the "Integral" / Red Violet palette. Only the theme color values
# differ between the two theme parts (dk1/lt1 are already identical in both),
# so the swap is expressed as updating each theme's ThemeColorScheme entries.

$p = $ppt.ActivePresentation

# --- Slide master's theme (ppt/theme/theme1.xml): adopt the "Office Theme" colors ---
$slideTheme = $p.SlideMaster.Theme.ThemeColorScheme
$slideTheme.Item(1).RGB  = 0         # dk1      000000
$slideTheme.Item(2).RGB  = 16777215  # lt1      FFFFFF
$slideTheme.Item(3).RGB  = 6968388   # dk2      44546A
$slideTheme.Item(4).RGB  = 15132391  # lt2      E7E6E6
$slideTheme.Item(5).RGB  = 13998939  # accent1  5B9BD5
$slideTheme.Item(6).RGB  = 3243501   # accent2  ED7D31
$slideTheme.Item(7).RGB  = 10855845  # accent3  A5A5A5
$slideTheme.Item(8).RGB  = 49407     # accent4  FFC000
$slideTheme.Item(9).RGB  = 12874308  # accent5  4472C4
$slideTheme.Item(10).RGB = 4697456   # accent6  70AD47
$slideTheme.Item(11).RGB = 12673797  # hlink    0563C1
$slideTheme.Item(12).RGB = 7491477   # folHlink 954F72

# NOTE: in this runtime, $p.NotesMaster.Theme.ThemeColorScheme resolves to the
# SAME underlying theme object as $p.SlideMaster.Theme.ThemeColorScheme (i.e.
# ppt/theme/theme1.xml) instead of the notes master's own ppt/theme/theme2.xml.
# Writing through it here would simply clobber the slide-master edit above
# instead of reaching theme2.xml, so it is intentionally not used.
